$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (shared-string) value updates in column G ("Correct")
$ws.Range("G2").Value = "318 (23)"
$ws.Range("G5").Value = "376 (34)"
$ws.Range("G8").Value = "352 (34)"

# Numeric updates in column E ("Error.Rate")
$ws.Range("E2").Value = 10.92436974789916
$ws.Range("E5").Value = 0.2652519893899204
$ws.Range("E8").Value = 4.3478260869565215

# Numeric updates in column I ("Miscellaneous")
$ws.Range("I2").Value = 167.0
$ws.Range("I3").Value = 147.0
$ws.Range("I4").Value = 229.0
$ws.Range("I5").Value = 192.0
$ws.Range("I6").Value = 196.0
$ws.Range("I7").Value = 149.0
$ws.Range("I8").Value = 199.0
$ws.Range("I9").Value = 183.0
